# fix name match happening for different id people
#
# The sheet had two sets of "-missingID" (E:M) and "-JKO" (V:AD) course
# columns that held duplicate data. The "-JKO" columns (V:AD) were the
# redundant/duplicate copy, so they are removed and the "-missingID"
# headers are relabeled to "-JKO" (the real source). A few emails that had
# been matched to the wrong person are also cleared/corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "-JKO" course columns (V:AD) entirely.
$ws.Range("V1:AD15").Delete()

# Relabel the remaining "-missingID" headers to "-JKO" (the correct source).
$ws.Range("E1").Value = "DHA Accommodations (1 hr)-JKO"
$ws.Range("F1").Value = "Leadership Training (4 hrs)-JKO"
$ws.Range("G1").Value = "MHS Customer Service (1 hr)-JKO"
$ws.Range("H1").Value = "Counterintelligence (1 hr)-JKO"
$ws.Range("I1").Value = "HIPAA Training (1 hr)-JKO"
$ws.Range("J1").Value = "Supervisor Safety Training (2 hrs)-JKO"
$ws.Range("K1").Value = "Employee Safety (1 hr)-JKO"
$ws.Range("L1").Value = "Violence Response (1 hr)-JKO"
$ws.Range("M1").Value = "RandomCourse-JKO"

# These emails had been matched to the wrong person by name alone; clear them.
$ws.Range("B2").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("B5").Value = ""

# Nick Fletcher's email was wrong; correct it.
$ws.Range("B4").Value = "n@gmail.com"
